# Generate Report for Archive
#
# 1) The status text "Ready for handoff" becomes "In Translation" everywhere
#    it is used (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
# 2) The "Status" columns (E/F on Overview, C on zh-cn/de-de) are narrowed
#    from their old, wider autofit width down to the new narrower width
#    that matches the shorter "In Translation" label.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update the status text -------------------------------------------------
$overview.Range("E2:F4").Value = "In Translation"
$zhcn.Range("C2:C4").Value = "In Translation"
$dede.Range("C2:C4").Value = "In Translation"

# --- Narrow the affected columns -------------------------------------------
$overview.Columns(5).ColumnWidth = 12.5   # column E
$overview.Columns(6).ColumnWidth = 12.5   # column F
$zhcn.Columns(3).ColumnWidth = 12.5       # column C
$dede.Columns(3).ColumnWidth = 12.5       # column C
